# Apply "special conditions expected results" update to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: description text changed from the single-condition note to the
#     "3 special conditions per layer" note, and the stray "For OED2..." note
#     in column I is removed.
$ws.Range("E5").Value = "3 special conditions per layer (7 policies)"
$ws.Range("I5").Value = ""

# --- Remove the other stray "For OED2, duplicates in location removed" notes
$ws.Range("I8").Value = ""
$ws.Range("I9").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("I13").Value = ""

# --- Rows 14-19: these test cases are now re-implemented ("yes") and done.
$ws.Range("G14:G19").Value = "yes"
$ws.Range("H14:H19").Value = "done"

# --- New row 20: sc15 test case
$ws.Range("A20").Value = "sc15"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = "Single hierarchy conditions One policy. One location with overlapping conditions"
$ws.Range("F20").Value = "input files"
$ws.Range("G20").Value = "no"
$ws.Range("H20").Value = "to do"

# --- Update the visible selection to match the saved workbook state
$ws.Range("H14:H19").Select()
